$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Elements")

$ws1.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/reference-with-code-and-period"
$ws1.Range("B3").Value = "8.0.0"
$ws1.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$ws1.Range("B9").Value = "LinuxForHealth Team"

$ws2.Range("AI2").Value = ""
$ws2.Range("J5").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/reference-code}`n"
$ws2.Range("J6").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/reference-period}`n"
